# Fix some text file issues:
# The "Text" column (C) for a few rows had embedded manual line breaks
# (wrapped at ~100 chars) left over from pasting. Replace those with the
# same text on a single logical line (line breaks collapsed to spaces),
# matching how the rest of the sheet's passages are stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C3").Value = "I agree that California's ""three strikes and you're out"" law will be a financial disaster for taxpayers who care about education and other vital services. But it's far from clear that the law can even be credited with a reduction in crime in California. While it's true that crime declined in California last year, crime also dropped nationwide."

$ws.Range("C4").Value = "Dorothy didn't know. She looked around her anxiously for some familiar landmark; but everything was strange. Between the branches of the many roads were green meadows and a few shrubs and trees, but she couldn't see the farm-house from which she had just come, or anything she had ever seen before, except the shaggy man and Toto."

$ws.Range("C5").Value = "Known as Rapa Nui to the island's inhabitants, Rongorongo is a writing system comprised of pictographs. It has been found carved into many oblong wooden tablets and other artifacts from the island's history. The art of writing was not known in any nearby islands and the script's mere existence is sufficient to confound anthropologists."

# The shorter, single-line text now wraps onto fewer visual lines, so the
# rows shrink from their old "4-line" autofit height down to "3 lines".
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45

# Leave the cursor where the author last left it.
$ws.Range("C16").Select() | Out-Null
